{"js": "// Worksheet was regenerated for the next day: the date heading advances\n// one day, and every \"A\u00d7B=C\" cell in the practice table is replaced with\n// a freshly generated three-digit-by-one-digit multiplication problem.\nconst replacements = [\n  [\"2024-11-09 Saturday\", \"2024-11-10 Sunday\"],\n  [\"520\u00d72=1040\", \"349\u00d76=2094\"],\n  [\"628\u00d76=3768\", \"410\u00d76=2460\"],\n  [\"155\u00d78=1240\", \"986\u00d79=8874\"],\n  [\"980\u00d73=2940\", \"919\u00d77=6433\"],\n  [\"712\u00d72=1424\", \"868\u00d74=3472\"],\n  [\"287\u00d72=574\", \"346\u00d76=2076\"],\n  [\"963\u00d77=6741\", \"947\u00d75=4735\"],\n  [\"529\u00d77=3703\", \"945\u00d74=3780\"],\n  [\"626\u00d76=3756\", \"412\u00d72=824\"],\n  [\"688\u00d77=4816\", \"547\u00d77=3829\"],\n  [\"727\u00d79=6543\", \"869\u00d75=4345\"],\n  [\"772\u00d75=3860\", \"820\u00d76=4920\"],\n  [\"431\u00d76=2586\", \"315\u00d75=1575\"],\n  [\"282\u00d74=1128\", \"335\u00d78=2680\"],\n  [\"425\u00d76=2550\", \"501\u00d76=3006\"],\n  [\"805\u00d77=5635\", \"112\u00d73=336\"],\n  [\"818\u00d77=5726\", \"695\u00d72=1390\"],\n  [\"249\u00d72=498\", \"263\u00d76=1578\"],\n  [\"869\u00d78=6952\", \"814\u00d77=5698\"],\n  [\"486\u00d75=2430\", \"410\u00d75=2050\"],\n  [\"853\u00d76=5118\", \"917\u00d76=5502\"],\n  [\"878\u00d74=3512\", \"444\u00d73=1332\"],\n  [\"532\u00d73=1596\", \"568\u00d74=2272\"],\n  [\"310\u00d79=2790\", \"396\u00d77=2772\"],\n  [\"862\u00d78=6896\", \"980\u00d74=3920\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all multiplication problems/answers in the table\n# to the next day's generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-11-09 Saturday', '2024-11-10 Sunday'),\n    @('520\u00d72=1040', '349\u00d76=2094'),\n    @('628\u00d76=3768', '410\u00d76=2460'),\n    @('155\u00d78=1240', '986\u00d79=8874'),\n    @('980\u00d73=2940', '919\u00d77=6433'),\n    @('712\u00d72=1424', '868\u00d74=3472'),\n    @('287\u00d72=574', '346\u00d76=2076'),\n    @('963\u00d77=6741', '947\u00d75=4735'),\n    @('529\u00d77=3703', '945\u00d74=3780'),\n    @('626\u00d76=3756', '412\u00d72=824'),\n    @('688\u00d77=4816', '547\u00d77=3829'),\n    @('727\u00d79=6543', '869\u00d75=4345'),\n    @('772\u00d75=3860', '820\u00d76=4920'),\n    @('431\u00d76=2586', '315\u00d75=1575'),\n    @('282\u00d74=1128', '335\u00d78=2680'),\n    @('425\u00d76=2550', '501\u00d76=3006'),\n    @('805\u00d77=5635', '112\u00d73=336'),\n    @('818\u00d77=5726', '695\u00d72=1390'),\n    @('249\u00d72=498', '263\u00d76=1578'),\n    @('869\u00d78=6952', '814\u00d77=5698'),\n    @('486\u00d75=2430', '410\u00d75=2050'),\n    @('853\u00d76=5118', '917\u00d76=5502'),\n    @('878\u00d74=3512', '444\u00d73=1332'),\n    @('532\u00d73=1596', '568\u00d74=2272'),\n    @('310\u00d79=2790', '396\u00d77=2772'),\n    @('862\u00d78=6896', '980\u00d74=3920'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n"}
